# Reloading code should be almost done. Just need to implement
# squareUpUsingProx(), centerOnIrArray(), pushButton().
#
# Spreadsheet side of the change: a new "Hall effect sensor" -> A6 entry
# was added to the "Usages" table (right after the two Proximity Sensor
# rows), pushing every subsequent row down by one. The "Arduino Mega
# Pins" sheet's A6 usage formula then flips from Free to Used
# automatically once that is in place.

$wb = $excel.ActiveWorkbook

$usages = $wb.Worksheets.Item("Usages")
$pins   = $wb.Worksheets.Item("Arduino Mega Pins")

# --- Usages sheet: insert the new Hall effect sensor row ------------------
# Row 42 currently starts the "IR Array" block; inserting here shifts it
# (and everything below) down to make room, inheriting the formatting of
# the row above (Proximity Sensor2), exactly like Excel's native
# Insert-row behavior.
$usages.Rows.Item(42).Insert()

$usages.Cells.Item(42, 1).Value = "Hall effect sensor"
$usages.Cells.Item(42, 3).Value = "A6"

# --- View/selection bookkeeping -------------------------------------------
# Touch the "Arduino Mega Pins" sheet first so its own scroll position can
# be set, then finish on "Usages" so it ends up the active tab/selection,
# matching the saved workbook state.
$pins.Activate()
$pinsWin = $excel.ActiveWindow
$pinsWin.ScrollRow = 47
$pinsWin.ScrollColumn = 1

$usages.Activate()
$usages.Range("E43").Select()
$usagesWin = $excel.ActiveWindow
$usagesWin.ScrollRow = 33
$usagesWin.ScrollColumn = 1
